$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.352.95"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "2.640.74"
$ws.Range("E3").Value = "  -3.44%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'550.07"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").Value = "'153.89"
$ws.Range("E6").Value = "  -4.46%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("E9").Value = "  -4.56%  "
$ws.Range("E10").Value = "  -4.33%  "
$ws.Range("D11").Value = "'5.46"
$ws.Range("E11").Value = "  -3.50%  "
$ws.Range("E12").Value = "  -4.95%  "
$ws.Range("D13").Value = "3.112.20"
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("D14").Value = "'25.65"
$ws.Range("E14").Value = "  -4.89%  "
$ws.Range("D15").Value = "62.269.18"
$ws.Range("E15").Value = "  -2.23%  "
$ws.Range("E16").Value = "  -4.27%  "
$ws.Range("D17").Value = "2.646.48"
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").Value = "'11.68"
$ws.Range("E18").Value = "  -5.97%  "
$ws.Range("E19").Value = "  -4.30%  "
$ws.Range("D20").Value = "'339.80"
$ws.Range("E20").Value = "  -4.38%  "
$ws.Range("E21").Value = "  -8.01%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'0.502"
$ws.Range("E23").Value = "  -3.73%  "
$ws.Range("D24").Value = "'62.64"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -4.77%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0834"
$ws.Range("E28").Value = "  -8.16%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "'1.35"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("E31").Value = "  -5.67%  "
$ws.Range("D32").Value = "'160.57"
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'4.73"
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("E35").Value = "  -3.56%  "
$ws.Range("D36").Value = "'19.19"
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("D37").Value = "'1.72"
$ws.Range("E37").Value = "  -5.03%  "
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("D39").Value = "'6.08"
$ws.Range("E39").Value = "  -3.56%  "
$ws.Range("D40").Value = "'0.906"
$ws.Range("E40").Value = "  -7.32%  "
$ws.Range("D41").Value = "'38.14"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("D42").Value = "'3.93"
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").Value = "'20.33"
$ws.Range("E44").Value = "  -6.57%  "
$ws.Range("D45").Value = "'0.607"
$ws.Range("E45").Value = "  -3.82%  "
$ws.Range("D46").Value = "'19.79"
$ws.Range("E46").Value = "  -6.56%  "
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").Value = "'0.0547"
$ws.Range("E48").Value = "  -6.82%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.0958"
$ws.Range("E49").Value = "  -4.19%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'127.60"
$ws.Range("E50").Value = "  -3.58%  "
$ws.Range("E51").Value = "  -5.81%  "
